$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 24:25 for the new HU#2_T6 / HU#2_T7 tasks
$ws.Rows("24:25").Insert()

$ws.Range("A24:D24").Style = "Normal"
$ws.Range("A24").Value = "HU#2_T6 Aplicar pre-commit en el proyecto."

$ws.Range("A25:D25").Style = "Normal"
$ws.Range("A25").Value = "HU#2_T7 Aplicar actions en el proyecto."

# New row 47 at the bottom of the table, thin left/right border
$ws.Rows("47:47").Insert()
$ws.Range("A47").Borders.Item(7).LineStyle = 1
$ws.Range("A47").Borders.Item(10).LineStyle = 1

# Update the sheet view (scrolled + selection moved)
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("A26").Select()
